# Student template: responsive improvement / cleanup of the "Tahun masuk"
# (K) column, fix two mis-typed Prodi values, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two rows had placeholder / typo'd text in the "Prodi" column (H) — both
# should just read "ES" like the rest of the sample rows.
$ws.Range("H6").Value = "ES"
$ws.Range("H7").Value = "ES"

# The "Tahun masuk" column (K) — header, its data validation-style cell
# comment, and all sample data — is being removed entirely.
$ws.Range("K1").Comment.Delete() | Out-Null
$ws.Columns.Item(11).Delete() | Out-Null

# Leave the selection where the author left it.
$ws.Range("K3").Select() | Out-Null
